$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("IDA")

# Row 14: Scenario text (B14) - remove "; Failure in Decryption"
$ws.Range("B14").Value = "Could not process request/Unknown error; Invalid Auth Request"

# Row 34: Message text (C34) - change decrypt message, highlight row 34 (B34:F34) yellow
$ws.Range("C34").Value = [char]8220 + "Unable to decrypt Request." + [char]8221

$ws.Range("B34:F34").Interior.Color = 65535

# Sheet view: scroll back to top, change active selection to J7
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("J7").Select()
